$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as published by the source feed.
$ws.Range("D2").Value = "67.690.73"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.329.53"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.42"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.74"
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "3.326.90"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.56"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "704.96"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "3.876.84"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.45"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "67.677.29"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "3.337.99"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.99"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.894"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  +3.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.93"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.23"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.51"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.18"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("E31").Value = "  +4.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "570.77"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.35"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "3.704.07"
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.02"
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.20"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0673"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.336"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  +5.22%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.71"
$ws.Range("E51").Value = "  +0.05%  "
